$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Prakash / 101 stay the same, only the scan time changes
$ws.Range("C2").Value = "21:40:49"

# Row 3: Mary / 102 / 21:41:06 (was Kolass / 103 / 17:49:08)
$ws.Range("A3").Value = "Mary"
$ws.Range("B3").Value = 102
$ws.Range("C3").Value = "21:41:06"

# Row 4: Prakash / 101 / 21:42:52 (was Mary / 102 / 17:49:28)
$ws.Range("A4").Value = "Prakash"
$ws.Range("B4").Value = 101
$ws.Range("C4").Value = "21:42:52"

# Row 5: Kolass / 103 / 21:43:04 (was Prakash / 101 / 17:50:30)
$ws.Range("A5").Value = "Kolass"
$ws.Range("B5").Value = 103
$ws.Range("C5").Value = "21:43:04"

# Row 6: Mary / "102" (text) / 21:43:04 (was Prakash / 101 / 18:01:04)
$ws.Range("A6").Value = "Mary"
$ws.Range("B6").Value = "'102"
$ws.Range("C6").Value = "21:43:04"

# The old rows 7-9 are no longer part of the rolling log, drop them so the
# sheet's used range shrinks back down to A1:C6
$ws.Range("A7:C9").ClearContents()
